$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.984.03"
$ws.Range("E2").Value = "  -4.83%  "
$ws.Range("D3").Value = "'2.484.48"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'535.63"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").Value = "'145.01"
$ws.Range("E6").Value = "  -6.31%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("D9").Value = "'2.511.21"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("D10").Value = "'0.0995"
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'5.40"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "'2.923.84"
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("D15").Value = "'23.93"
$ws.Range("E15").Value = "  -6.53%  "
$ws.Range("D16").Value = "'58.906.46"
$ws.Range("E16").Value = "  -4.86%  "
$ws.Range("E17").Value = "  -3.98%  "
$ws.Range("D18").Value = "'2.511.54"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("D19").Value = "'11.26"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("E20").Value = "  -5.87%  "
$ws.Range("D21").Value = "'323.49"
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'5.74"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").Value = "'61.36"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "'0.439"
$ws.Range("E25").Value = "  -10.89%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "'2.617.83"
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.161"
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("D29").Value = "'7.76"
$ws.Range("E29").Value = "  -5.50%  "
$ws.Range("D30").Value = "'6.83"
$ws.Range("E30").Value = "  -6.19%  "
$ws.Range("D31").Value = "'0.0₃0777"
$ws.Range("E31").Value = "  -7.29%  "
$ws.Range("E32").Value = "  -7.68%  "
$ws.Range("E33").Value = "  -5.45%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'158.73"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").Value = "'1.44"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").Value = "'18.50"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("D38").Value = "'4.43"
$ws.Range("E38").Value = "  -9.44%  "
$ws.Range("E39").Value = "  -10.05%  "
$ws.Range("D40").Value = "'5.87"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").Value = "'307.19"
$ws.Range("D42").Value = "'36.74"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").Value = "'3.67"
$ws.Range("E43").Value = "  -7.09%  "
$ws.Range("D44").Value = "'0.820"
$ws.Range("E44").Value = "  -9.56%  "
$ws.Range("D45").Value = "'0.994"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "'0.594"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").Value = "'10.78"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "'124.10"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "'0.0929"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("D50").Value = "'18.60"
$ws.Range("E50").Value = "  -4.85%  "
$ws.Range("D51").Value = "'0.0515"
$ws.Range("E51").Value = "  -5.99%  "

# Reset number format on the Price column so the forced-text apostrophe
# does not leave a stray text-format style applied to the cells.
$ws.Range("D2:D51").Style = "Normal"

